$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 4061.5454
$ws.Cells.Item(32, 9).Value = 4378
$ws.Cells.Item(32, 10).Value = 3797.8333
$ws.Cells.Item(32, 11).Value = 4378
$ws.Cells.Item(32, 12).Value = 3797.8333
$ws.Cells.Item(32, 13).Value = -4052
$ws.Cells.Item(32, 14).Value = -4449.8333
$ws.Cells.Item(132, 8).Value = 3806.1853
$ws.Cells.Item(132, 9).Value = 3921.8076
$ws.Cells.Item(132, 11).Value = 11765.4228
$ws.Cells.Item(132, 13).Value = -9235.4228
$ws.Cells.Item(137, 8).Value = 17243134
$ws.Cells.Item(137, 9).Value = 1020.7742
$ws.Cells.Item(137, 10).Value = 37039630
$ws.Cells.Item(137, 11).Value = 3062.3226
$ws.Cells.Item(137, 12).Value = 111118890
$ws.Cells.Item(137, 13).Value = -512.3226
$ws.Cells.Item(137, 14).Value = -111123990
$ws.Cells.Item(141, 8).Value = 1101.2858
$ws.Cells.Item(141, 9).Value = 1101.2858
$ws.Cells.Item(141, 10).Value = 0
$ws.Cells.Item(141, 11).Value = 3303.8574
$ws.Cells.Item(141, 12).Value = 0
$ws.Cells.Item(141, 13).Value = 1876.1426
$ws.Cells.Item(141, 14).Value = ""

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 2605.018
$ws.Cells.Item(61, 9).Value = 2215.9143
$ws.Cells.Item(61, 10).Value = 3285.95
$ws.Cells.Item(61, 11).Value = 2215.9143
$ws.Cells.Item(61, 12).Value = 3285.95
$ws.Cells.Item(61, 13).Value = -2003.9143
$ws.Cells.Item(61, 14).Value = -3709.95
$ws.Cells.Item(74, 8).Value = 5146.3516
$ws.Cells.Item(74, 9).Value = 1148.591
$ws.Cells.Item(74, 10).Value = 11009.733
$ws.Cells.Item(74, 11).Value = 1148.591
$ws.Cells.Item(74, 12).Value = 11009.733
$ws.Cells.Item(74, 13).Value = -274.5909999999999
$ws.Cells.Item(74, 14).Value = -12757.733
$ws.Cells.Item(77, 8).Value = 5146.3516
$ws.Cells.Item(77, 9).Value = 1148.591
$ws.Cells.Item(77, 10).Value = 11009.733
$ws.Cells.Item(77, 11).Value = 5742.955
$ws.Cells.Item(77, 12).Value = 55048.665
$ws.Cells.Item(77, 13).Value = -1374.955
$ws.Cells.Item(77, 14).Value = -63784.665
$ws.Cells.Item(136, 8).Value = 2605.018
$ws.Cells.Item(136, 9).Value = 2215.9143
$ws.Cells.Item(136, 10).Value = 3285.95
$ws.Cells.Item(136, 11).Value = 6647.742899999999
$ws.Cells.Item(136, 12).Value = 9857.849999999999
$ws.Cells.Item(136, 13).Value = -4097.742899999999
$ws.Cells.Item(136, 14).Value = -14957.85

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 194.27586
$ws.Cells.Item(80, 9).Value = 38.75
$ws.Cells.Item(80, 10).Value = 253.5238
$ws.Cells.Item(80, 11).Value = 38.75
$ws.Cells.Item(80, 12).Value = 253.5238
$ws.Cells.Item(80, 13).Value = 959.25
$ws.Cells.Item(80, 14).Value = -2249.5238
$ws.Cells.Item(83, 8).Value = 194.27586
$ws.Cells.Item(83, 9).Value = 38.75
$ws.Cells.Item(83, 10).Value = 253.5238
$ws.Cells.Item(83, 11).Value = 193.75
$ws.Cells.Item(83, 12).Value = 1267.619
$ws.Cells.Item(83, 13).Value = 4798.25
$ws.Cells.Item(83, 14).Value = -11251.619

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3270545.5
$ws.Cells.Item(31, 9).Value = 1213.8948
$ws.Cells.Item(31, 10).Value = 5211711.5
$ws.Cells.Item(31, 11).Value = 1213.8948
$ws.Cells.Item(31, 12).Value = 5211711.5
$ws.Cells.Item(31, 13).Value = -918.8948
$ws.Cells.Item(31, 14).Value = -5212301.5
$ws.Cells.Item(34, 8).Value = 3270545.5
$ws.Cells.Item(34, 9).Value = 1213.8948
$ws.Cells.Item(34, 10).Value = 5211711.5
$ws.Cells.Item(34, 11).Value = 1213.8948
$ws.Cells.Item(34, 12).Value = 5211711.5
$ws.Cells.Item(34, 13).Value = -1011.8948
$ws.Cells.Item(34, 14).Value = -5212115.5
$ws.Cells.Item(58, 8).Value = 2985.2837
$ws.Cells.Item(58, 9).Value = 1365.5385
$ws.Cells.Item(58, 10).Value = 4012.439
$ws.Cells.Item(58, 11).Value = 1365.5385
$ws.Cells.Item(58, 12).Value = 4012.439
$ws.Cells.Item(58, 13).Value = -1162.5385
$ws.Cells.Item(58, 14).Value = -4418.439
$ws.Cells.Item(69, 8).Value = 5930.3335
$ws.Cells.Item(69, 9).Value = 4796.625
$ws.Cells.Item(69, 10).Value = 15000
$ws.Cells.Item(69, 11).Value = 4796.625
$ws.Cells.Item(69, 12).Value = 15000
$ws.Cells.Item(69, 13).Value = -4047.625
$ws.Cells.Item(69, 14).Value = -16498
$ws.Cells.Item(72, 8).Value = 5930.3335
$ws.Cells.Item(72, 9).Value = 4796.625
$ws.Cells.Item(72, 10).Value = 15000
$ws.Cells.Item(72, 11).Value = 14389.875
$ws.Cells.Item(72, 12).Value = 45000
$ws.Cells.Item(72, 13).Value = -10645.875
$ws.Cells.Item(72, 14).Value = -52488
$ws.Cells.Item(132, 8).Value = 1596.7344
$ws.Cells.Item(132, 9).Value = 1482.6
$ws.Cells.Item(132, 10).Value = 2294.2222
$ws.Cells.Item(132, 11).Value = 4447.799999999999
$ws.Cells.Item(132, 12).Value = 6882.6666
$ws.Cells.Item(132, 13).Value = -1917.799999999999
$ws.Cells.Item(132, 14).Value = -11942.6666
$ws.Cells.Item(136, 8).Value = 2985.2837
$ws.Cells.Item(136, 9).Value = 1365.5385
$ws.Cells.Item(136, 10).Value = 4012.439
$ws.Cells.Item(136, 11).Value = 4096.6155
$ws.Cells.Item(136, 12).Value = 12037.317
$ws.Cells.Item(136, 13).Value = -1546.6155
$ws.Cells.Item(136, 14).Value = -17137.317

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 879.5862
$ws.Cells.Item(5, 9).Value = 674.5263
$ws.Cells.Item(5, 10).Value = 979.4872
$ws.Cells.Item(5, 11).Value = 2023.5789
$ws.Cells.Item(5, 12).Value = 2938.4616
$ws.Cells.Item(5, 13).Value = -1911.5789
$ws.Cells.Item(5, 14).Value = -3162.4616
$ws.Cells.Item(103, 8).Value = 2267439.2
$ws.Cells.Item(103, 10).Value = 962.55554
$ws.Cells.Item(103, 12).Value = 2887.66662
$ws.Cells.Item(103, 14).Value = -4645.66662
$ws.Cells.Item(131, 8).Value = 2082.8674
$ws.Cells.Item(131, 10).Value = 1413.2784
$ws.Cells.Item(131, 12).Value = 4239.8352
$ws.Cells.Item(131, 14).Value = -14319.8352
$ws.Cells.Item(135, 8).Value = 879.5862
$ws.Cells.Item(135, 9).Value = 674.5263
$ws.Cells.Item(135, 10).Value = 979.4872
$ws.Cells.Item(135, 11).Value = 6070.736699999999
$ws.Cells.Item(135, 12).Value = 8815.3848
$ws.Cells.Item(135, 13).Value = -3535.736699999999
$ws.Cells.Item(135, 14).Value = -13885.3848

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value = 477.95
$ws.Cells.Item(107, 9).Value = 311.5
$ws.Cells.Item(107, 11).Value = 311.5
$ws.Cells.Item(107, 13).Value = 1608.5
$ws.Cells.Item(126, 8).Value = 2170.8333
$ws.Cells.Item(126, 9).Value = 1891.6666
$ws.Cells.Item(126, 10).Value = 2450
$ws.Cells.Item(126, 11).Value = 5674.9998
$ws.Cells.Item(126, 12).Value = 7350
$ws.Cells.Item(126, 13).Value = -3204.9998
$ws.Cells.Item(126, 14).Value = -12290

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 4025.7917
$ws.Cells.Item(132, 9).Value = 4589.5835
$ws.Cells.Item(132, 10).Value = 2334.4167
$ws.Cells.Item(132, 11).Value = 13768.7505
$ws.Cells.Item(132, 12).Value = 7003.250100000001
$ws.Cells.Item(132, 13).Value = -11238.7505
$ws.Cells.Item(132, 14).Value = -12063.2501

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 14545199
$ws.Cells.Item(132, 9).Value = 11891085
$ws.Cells.Item(132, 10).Value = 28574086
$ws.Cells.Item(132, 11).Value = 35673255
$ws.Cells.Item(132, 12).Value = 85722258
$ws.Cells.Item(132, 13).Value = -35670725
$ws.Cells.Item(132, 14).Value = -85727318
